$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.087.46'
$ws.Range("D2").ClearFormats()

$ws.Range("E2").Value = '  -0.41%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.405.19'
$ws.Range("D3").ClearFormats()

$ws.Range("E3").Value = '  -0.86%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.53'
$ws.Range("D5").ClearFormats()

$ws.Range("E5").Value = '  -0.92%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.25'
$ws.Range("D6").ClearFormats()

$ws.Range("E6").Value = '  -1.66%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.406.51'
$ws.Range("D7").ClearFormats()

$ws.Range("E7").Value = '  -0.90%  '

$ws.Range("E9").Value = '  +0.21%  '

$ws.Range("E10").Value = '  -1.11%  '

$ws.Range("E11").Value = '  -0.61%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.395'
$ws.Range("D12").ClearFormats()

$ws.Range("E12").Value = '  +2.45%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.984.76'
$ws.Range("D13").ClearFormats()

$ws.Range("E13").Value = '  -1.00%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.52'
$ws.Range("D14").ClearFormats()

$ws.Range("E14").Value = '  +1.70%  '

$ws.Range("E15").Value = '  +1.93%  '

$ws.Range("E16").Value = '  -0.30%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.393.70'
$ws.Range("D17").ClearFormats()

$ws.Range("E17").Value = '  -2.63%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.150.42'
$ws.Range("D18").ClearFormats()

$ws.Range("E18").Value = '  -0.60%  '

$ws.Range("E19").Value = '  -0.60%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.00'
$ws.Range("D20").ClearFormats()

$ws.Range("E20").Value = '  -1.19%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.11'
$ws.Range("D21").ClearFormats()

$ws.Range("E21").Value = '  -4.02%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '387.02'
$ws.Range("D22").ClearFormats()

$ws.Range("E22").Value = '  -1.86%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.21'
$ws.Range("D24").ClearFormats()

$ws.Range("E24").Value = '  +1.64%  '

$ws.Range("E25").Value = '  +0.40%  '

$ws.Range("E26").Value = '  -3.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.538.18'
$ws.Range("D27").ClearFormats()

$ws.Range("E27").Value = '  -1.31%  '

$ws.Range("E28").Value = '  +0.80%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").ClearFormats()

$ws.Range("E29").Value = '  -0.11%  '

$ws.Range("E30").Value = '  -1.77%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.01'
$ws.Range("D31").ClearFormats()

$ws.Range("E31").Value = '  -1.53%  '

$ws.Range("E32").Value = '  -0.81%  '

$ws.Range("E33").Value = '  -3.30%  '

$ws.Range("E34").Value = '  -0.06%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.69'
$ws.Range("D35").ClearFormats()

$ws.Range("E35").Value = '  -1.21%  '

$ws.Range("E36").Value = '  +0.53%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '166.80'
$ws.Range("D37").ClearFormats()

$ws.Range("E37").Value = '  -0.40%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.435.71'
$ws.Range("D38").ClearFormats()

$ws.Range("E38").Value = '  -0.94%  '

$ws.Range("E39").Value = '  -1.49%  '

$ws.Range("E40").Value = '  -3.72%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '28.55'
$ws.Range("D41").ClearFormats()

$ws.Range("E41").Value = '  +4.45%  '

$ws.Range("E42").Value = '  +0.28%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.781'
$ws.Range("D43").ClearFormats()

$ws.Range("E43").Value = '  -2.40%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("D44").ClearFormats()

$ws.Range("E44").Value = '  -0.12%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.22'
$ws.Range("D45").ClearFormats()

$ws.Range("E45").Value = '  +0.42%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.45'
$ws.Range("D46").ClearFormats()

$ws.Range("E46").Value = '  -0.52%  '

$ws.Range("E47").Value = '  -2.65%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.13'
$ws.Range("D48").ClearFormats()

$ws.Range("E48").Value = '  -1.73%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.504.52'
$ws.Range("D49").ClearFormats()

$ws.Range("E49").Value = '  -3.47%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.38'
$ws.Range("D50").ClearFormats()

$ws.Range("E50").Value = '  +1.77%  '

$ws.Range("E51").Value = '  -0.91%  '
